# Applies the "data through 2018 and edited methods" update:
#  - ColumnHeaders: rename the picoeuk abundance attribute to a concentration
#    attribute with an updated definition.
#  - Personnel: add Robert J Olson (creator), Alexi Shalapyonok (technician)
#    and E. Taylor Crockford (technician); move Bethany Fowler down and
#    change her role to metadataProvider.
#  - Update sheet selections / active sheet to match the saved UI state.

$wb = $excel.ActiveWorkbook

$wsHeaders = $wb.Worksheets.Item("ColumnHeaders")
$wsPersonnel = $wb.Worksheets.Item("Personnel")
$wsKeywords = $wb.Worksheets.Item("Keywords")

# ---------------------------------------------------------------------------
# Personnel sheet
# ---------------------------------------------------------------------------
# Row 2 stays Heidi M Sosik (unchanged content, just re-set for clarity).
$wsPersonnel.Cells.Item(2,1).Value = "Heidi"
$wsPersonnel.Cells.Item(2,2).Value = "M"
$wsPersonnel.Cells.Item(2,3).Value = "Sosik"
$wsPersonnel.Cells.Item(2,4).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(2,5).Value = "hsosik@whoi.edu"
$wsPersonnel.Cells.Item(2,6).Value = "0000-0002-4591-2842"
$wsPersonnel.Cells.Item(2,7).Value = "creator"
$wsPersonnel.Cells.Item(2,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(2,9).Value = "NSF"
$wsPersonnel.Cells.Item(2,10).Value = "OCE-1655686"

# Row 3 (new): Robert J Olson, WHOI, creator.
$wsPersonnel.Cells.Item(3,1).Value = "Robert"
$wsPersonnel.Cells.Item(3,3).Value = "Olson"
$wsPersonnel.Cells.Item(3,2).Value = "J"
$wsPersonnel.Cells.Item(3,4).Value = "Woods Hole Oceanographic Institution"
$wsPersonnel.Cells.Item(3,6).Value = ""
$wsPersonnel.Cells.Item(3,7).Value = "creator"
$wsPersonnel.Cells.Item(3,8).Value = ""
$wsPersonnel.Cells.Item(3,9).Value = ""
$wsPersonnel.Cells.Item(3,10).Value = ""

# Row 4 stays the NES-LTER Information Manager contact (unchanged content).
$wsPersonnel.Cells.Item(4,1).Value = "NES-LTER Information Manager"
$wsPersonnel.Cells.Item(4,4).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(4,5).Value = "lter-nes-info@whoi.edu"
$wsPersonnel.Cells.Item(4,7).Value = "contact"
$wsPersonnel.Cells.Item(4,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(4,9).Value = "NSF"
$wsPersonnel.Cells.Item(4,10).Value = "OCE-1655686"

# Row 6 (new): E. Taylor Crockford, technician.
$wsPersonnel.Cells.Item(6,1).Value = "E. Taylor"
$wsPersonnel.Cells.Item(6,3).Value = "Crockford"
$wsPersonnel.Cells.Item(6,4).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(6,5).Value = "ecrockford@whoi.edu"
$wsPersonnel.Cells.Item(6,7).Value = "technician"
$wsPersonnel.Cells.Item(6,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(6,9).Value = "NSF"
$wsPersonnel.Cells.Item(6,10).Value = "OCE-1655686"

# Row 5 (new): Alexi Shalapyonok, technician.
$wsPersonnel.Cells.Item(5,1).Value = "Alexi"
$wsPersonnel.Cells.Item(5,3).Value = "Shalapyonok"
$wsPersonnel.Cells.Item(5,4).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(5,5).Value = "ashalapyonok@whoi.edu"
$wsPersonnel.Cells.Item(5,7).Value = "technician"
$wsPersonnel.Cells.Item(5,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(5,9).Value = "NSF"
$wsPersonnel.Cells.Item(5,10).Value = "OCE-1655686"

# Row 7 (was row 2): Bethany Fowler, now metadataProvider instead of creator.
$wsPersonnel.Cells.Item(7,1).Value = "Bethany"
$wsPersonnel.Cells.Item(7,3).Value = "Fowler"
$wsPersonnel.Cells.Item(7,4).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(7,5).Value = "bfowler@whoi.edu"
$wsPersonnel.Cells.Item(7,7).Value = "metadataProvider"
$wsPersonnel.Cells.Item(7,8).Value = "Northeast U.S. Shelf LTER"
$wsPersonnel.Cells.Item(7,9).Value = "NSF"
$wsPersonnel.Cells.Item(7,10).Value = "OCE-1655686"

# ---------------------------------------------------------------------------
# ColumnHeaders sheet: rename abundance_picoeuk -> concentration_picoeuk
# ---------------------------------------------------------------------------
$wsHeaders.Cells.Item(2,1).Value = "concentration_picoeuk"
$wsHeaders.Cells.Item(2,2).Value = "abundance per unit volume of eukaryote picophytoplankton http://vocab.nerc.ac.uk/collection/P01/current/SDBIOL01/ http://vocab.nerc.ac.uk/collection/F02/current/F0200004/"
$wsHeaders.Cells.Item(2,3).Value = "numeric"
$wsHeaders.Cells.Item(2,4).Value = "numberPerMilliliter"
$wsHeaders.Cells.Item(2,6).Value = "NaN"
$wsHeaders.Cells.Item(2,7).Value = "Missing value"

# Robert Olson's email is filled in last, after the other new strings above.
$wsPersonnel.Cells.Item(3,5).Value = "rolson@whoi.edu"

# ---------------------------------------------------------------------------
# View state: ColumnHeaders selection -> B2, Personnel becomes the active
# tab with selection E5, Keywords keeps its F4 selection but is no longer
# the active tab.
# ---------------------------------------------------------------------------
$wsHeaders.Activate()
$wsHeaders.Range("B2").Select()

$wsKeywords.Activate()
$wsKeywords.Range("F4").Select()

$wsPersonnel.Activate()
$wsPersonnel.Range("E5").Select()
